$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4113.531486898866
$ws.Range("C3").Value = 4113.531486898866
$ws.Range("C4").Value = 4113.531486898866
$ws.Range("C5").Value = 4108.251268283465
$ws.Range("C6").Value = 4108.251268283465
$ws.Range("C7").Value = 4108.251268283465
$ws.Range("C8").Value = 4036.781481655154
$ws.Range("C9").Value = 4036.781481655154
$ws.Range("C10").Value = 4036.781481655154
$ws.Range("C11").Value = 3929.144662485524
$ws.Range("C12").Value = 3929.144662485524
